# "Update newly subscripted variables in CID file"
#
# The "Boolean" tab (3rd worksheet) lists InputData CSV pathnames that are
# checked for boolean-constrained values. Two of those pathnames
# (trans/BVTQaZ/BVTQaZ.csv and trans/VTQaZ/VTQaZ.csv) have been split up
# ("subscripted") into six per-mode files each: LDVs, HDVs, aircraft, rail,
# ships, motorbikes. This script removes the old single-file rows and
# inserts the six new rows in their place, for both BVTQaZ and VTQaZ.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Boolean")

# --- trans/BVTQaZ/BVTQaZ.csv (currently row 17) -> 6 subscripted rows ---
# Insert 5 blank rows below row 17 so rows 17-22 are available, then fill
# them with the new subscripted pathnames (alphabetical, same order as the
# sharedStrings additions in the target workbook).
$ws.Rows("18:22").Insert()

$ws.Range("A17").Value = "trans/BVTQaZ/BVTQaZ-LDVs.csv"
$ws.Range("A18").Value = "trans/BVTQaZ/BVTQaZ-HDVs.csv"
$ws.Range("A19").Value = "trans/BVTQaZ/BVTQaZ-aircraft.csv"
$ws.Range("A20").Value = "trans/BVTQaZ/BVTQaZ-rail.csv"
$ws.Range("A21").Value = "trans/BVTQaZ/BVTQaZ-ships.csv"
$ws.Range("A22").Value = "trans/BVTQaZ/BVTQaZ-motorbikes.csv"

# --- trans/VTQaZ/VTQaZ.csv (was row 21, now shifted down to row 26) ---
$ws.Rows("27:31").Insert()

$ws.Range("A26").Value = "trans/VTQaZ/VTQaZ-LDVs.csv"
$ws.Range("A27").Value = "trans/VTQaZ/VTQaZ-HDVs.csv"
$ws.Range("A28").Value = "trans/VTQaZ/VTQaZ-aircraft.csv"
$ws.Range("A29").Value = "trans/VTQaZ/VTQaZ-rail.csv"
$ws.Range("A30").Value = "trans/VTQaZ/VTQaZ-ships.csv"
$ws.Range("A31").Value = "trans/VTQaZ/VTQaZ-motorbikes.csv"

# A handful of trailing blank (but row-formatted) rows follow the last
# entry in the refreshed workbook.
$ws.Rows("33:38").RowHeight = 15

# --- View/selection bookkeeping that Excel records when a user works on
# the workbook (matches the author's save: they ended on the "Boolean"
# sheet at A32, with "Integer" parked at A13, and "About" as the active
# tab when the file was reopened/saved). ---
$ws.Range("A32").Select()

$wsInteger = $wb.Worksheets.Item("Integer")
$wsInteger.Range("A13").Select()

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
